# Finalização da V1 do programa da shalom
# Adds a "Link" column (F) with YouTube links for each song in the
# "musicas" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell F1 ("Link") -------------------------------------------
# Match the formatting already used by the other header cells (A1:E1):
# copy the format from E1 (grey fill, border, centered) onto F1, then set
# its text.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Link"

# --- Data cells F2:F18 ---------------------------------------------------
# These use a simple bordered style (no special alignment), matching the
# rest of the data rows' thin-border look.
$linkRange = $ws.Range("F2:F18")
$linkRange.Borders.ColorIndex = 1
$linkRange.Borders.LineStyle = 1

$links = @(
  "https://www.youtube.com/watch?v=7DH_tKN_n-g",
  "https://www.youtube.com/watch?v=4GC0uxYbJ-M",
  "https://www.youtube.com/watch?v=x1h5h1VWN6Y",
  "https://www.youtube.com/watch?v=K9wXc0FWITM",
  "https://www.youtube.com/watch?v=nv-T2_JPKZA",
  "https://www.youtube.com/watch?v=JheqX_w3m08",
  "https://www.youtube.com/watch?v=BqzYurSgAtk",
  "https://www.youtube.com/watch?v=b0lQYSjyMfM",
  "https://www.youtube.com/watch?v=60CuZyzGf5U",
  "https://www.youtube.com/watch?v=1DzK7Wm3IcE",
  "https://www.youtube.com/watch?v=oAyIGD3Ek7g",
  "https://www.youtube.com/watch?v=3JUS_ueGjnA",
  "https://www.youtube.com/watch?v=KhVSYlfiL84",
  "https://www.youtube.com/watch?v=mZ9yZYo9Mmk",
  "https://www.youtube.com/watch?v=cbAu_85RRtc",
  "https://www.youtube.com/watch?v=QbnmpJo3DiI",
  "https://www.youtube.com/watch?v=5-srhflJ-kg"
)

for ($i = 0; $i -lt $links.Length; $i++) {
  $row = $i + 2
  $cellRef = "F" + $row
  $ws.Range($cellRef).Value = $links[$i]
}

# Size the new column to fit its (long) URL contents.
$ws.Columns.Item(6).AutoFit()

# Leave the selection where the author ended up after entering the data.
$ws.Range("E26").Select() | Out-Null
